# Fruta / hortaliza, semanal
# A new daily-price record was inserted as row 250 of the "Mandarina" sheet,
# pushing the previously existing rows 250-335 down to 251-336.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 250 (shifts rows 250:335 -> 251:336)
$ws.Rows(250).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = 44809
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = "Fruta"
$ws.Range("G250").Value = 100102
$ws.Range("H250").Value = "Cítricos"
$ws.Range("I250").Value = 100102004
$ws.Range("J250").Value = "Mandarina"
$ws.Range("K250").Value = "Murcott"
$ws.Range("L250").Value = "Primera"
$ws.Range("M250").Value = 300
$ws.Range("N250").Value = 7000
$ws.Range("O250").Value = 7000
$ws.Range("P250").Value = 7000
$ws.Range("Q250").Value = "`$/bandeja 10 kilos"
$ws.Range("R250").Value = "Provincia de Limarí"
$ws.Range("S250").Value = 700
$ws.Range("T250").Value = 10
